# Apply the changes described in the commit:
# "got supporting entity visualization now differentiates between
#  non-numerical and numerical answers."
#
# 1. "SemScores Analysis" sheet: add a new summary row 26 ("Overall"/"ALL")
# 2. "Worst Performers" sheet: T6 and T7 (hit_1 column) flip from False to True
# 3. "Hits_1 Semscores" sheet: updated hit/not_hit mean & variance values

$wb = $excel.ActiveWorkbook

# --- 1. "SemScores Analysis" sheet: append row 26 ---
$wsAnalysis = $wb.Worksheets.Item("SemScores Analysis")
$wsAnalysis.Range("A26").Value = "Overall"
$wsAnalysis.Range("B26").Value = "ALL"
$wsAnalysis.Range("C26").Value = 0.5844565208755611
$wsAnalysis.Range("D26").Value = 0.08162608204393407

# --- 2. "Worst Performers" sheet: hit_1 flag flips for rows 6 and 7 ---
# (values are stored as *text* "True"/"False", not Excel booleans, so a plain
#  Value assignment of the string "True" would be auto-coerced to a boolean;
#  instead write a formula producing the text and paste back as a literal value)
$wsWorst = $wb.Worksheets.Item("Worst Performers")

$cellT6 = $wsWorst.Range("T6")
$cellT6.Formula = "=""True"""
$cellT6.Copy()
$cellT6.PasteSpecial(-4163)

$cellT7 = $wsWorst.Range("T7")
$cellT7.Formula = "=""True"""
$cellT7.Copy()
$cellT7.PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- 3. "Hits_1 Semscores" sheet: updated mean/variance values ---
$wsHits = $wb.Worksheets.Item("Hits_1 Semscores")
$wsHits.Range("B2").Value = 0.6406397391691465
$wsHits.Range("C2").Value = 0.06235470815222292
$wsHits.Range("B3").Value = 0.5510683306416638
$wsHits.Range("C3").Value = 0.08405391080859866
